# Regenerate save_data: update column G ("K", strikeout count) values
# per the recomputed std/mean and s_vals.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 0
    3  = 3
    4  = 0
    5  = 2
    6  = 0
    7  = 0
    8  = 2
    9  = 1
    10 = 1
    11 = 1
    12 = 1
    13 = 0
    14 = 3
    15 = 2
    16 = 2
    17 = 1
    18 = 2
    19 = 1
    20 = 2
    21 = 0
    22 = 1
    23 = 2
    24 = 2
    25 = 0
    26 = 0
    27 = 0
    28 = 2
    29 = 1
    30 = 2
    31 = 3
    32 = 0
    33 = 1
    34 = 0
    35 = 0
    36 = 1
    37 = 2
    38 = 1
    39 = 0
    40 = 2
    41 = 3
    42 = 0
    43 = 1
    44 = 1
    45 = 2
    46 = 2
    47 = 2
    48 = 1
    49 = 3
    50 = 1
    51 = 3
    52 = 1
    53 = 0
    54 = 1
    55 = 1
    56 = 2
    57 = 0
    58 = 0
    59 = 0
    60 = 2
    61 = 3
    62 = 0
    63 = 1
    64 = 0
    65 = 1
    66 = 1
    67 = 3
    68 = 0
    69 = 1
    70 = 1
    71 = 0
    72 = 1
    73 = 2
    74 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
